$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "worker 13, worker 15"
$ws.Range("C2").Value = "worker 13, worker 15"
$ws.Range("D2").Value = "worker 3, worker 13, worker 15"
$ws.Range("E2").Value = "worker 13, worker 15"
$ws.Range("F2").Value = "worker 13, worker 15"

$ws.Range("B3").Value = "worker 4, worker 9"
$ws.Range("C3").Value = "worker 2, worker 6"
$ws.Range("D3").Value = "worker 3, worker 9"
$ws.Range("E3").Value = "worker 2, worker 6"
$ws.Range("F3").Value = "worker 3, worker 6"

$ws.Range("B4").Value = "worker 8, worker 9"
$ws.Range("C4").Value = "worker 2, worker 12"
$ws.Range("D4").Value = "worker 8, worker 9"
$ws.Range("E4").Value = "worker 2, worker 16"
$ws.Range("F4").Value = "worker 6, worker 9"

$ws.Range("B5").Value = "worker 13, worker 14"
$ws.Range("D5").Value = "worker 8, worker 15"
$ws.Range("F5").Value = "worker 5, worker 11"

$ws.Range("B6").Value = "worker 4, worker 7, worker 11, worker 14, worker 16"
$ws.Range("C6").Value = "worker 2, worker 3, worker 7, worker 8, worker 11"
$ws.Range("D6").Value = "worker 3, worker 4, worker 7, worker 14, worker 15"
$ws.Range("E6").Value = "worker 2, worker 3, worker 8, worker 11, worker 14"
$ws.Range("F6").Value = "worker 1, worker 3, worker 4, worker 6, worker 9"

$ws.Range("B7").Value = "worker 1, worker 7, worker 10, worker 11, worker 16"
$ws.Range("C7").Value = "worker 1, worker 6, worker 7, worker 11, worker 13"
$ws.Range("D7").Value = "worker 1, worker 4, worker 7, worker 10, worker 11"
$ws.Range("E7").Value = "worker 1, worker 6, worker 10, worker 11, worker 14"
$ws.Range("F7").Value = "worker 1, worker 5, worker 6, worker 9, worker 10"

$ws.Range("B8").Value = "worker 3, worker 4, worker 11, worker 16"
$ws.Range("C8").Value = "worker 4, worker 9, worker 12, worker 13"
$ws.Range("D8").Value = "worker 1, worker 3, worker 7, worker 12"
$ws.Range("E8").Value = "worker 1, worker 4, worker 9, worker 11"
$ws.Range("F8").Value = "worker 1, worker 10, worker 11, worker 12"

$ws.Range("B9").Value = "worker 2, worker 4, worker 7, worker 12"
$ws.Range("C9").Value = "worker 5, worker 8, worker 9, worker 11"
$ws.Range("D9").Value = "worker 1, worker 2, worker 3, worker 8"
$ws.Range("E9").Value = "worker 4, worker 9, worker 11, worker 14"
$ws.Range("F9").Value = "worker 1, worker 2, worker 4, worker 8"

$ws.Range("B10").Value = "worker 2, worker 4, worker 14"
$ws.Range("C10").Value = "worker 3, worker 11, worker 14"
$ws.Range("D10").Value = "worker 2, worker 9, worker 11"
$ws.Range("E10").Value = "worker 3, worker 8, worker 14"
$ws.Range("F10").Value = "worker 1, worker 2, worker 8"

$wb.Save()
